# sprint_project.xlsx — "Add files via upload" edit
#
# Net effect of the commit (per the OOXML diff): the Status cell for the
# first task row (E2) gets filled in with "Done", centered both
# horizontally and vertically, matching the look of the other filled-in
# cells in that column's header/border style.
#
# (The diff also shows a lot of incidental churn — fileVersion/rupBuild,
# absPath, revisionPtr/coauth versions, window geometry, calcId, the
# default theme font rendering as Calibri instead of Arial, and the
# resulting default-row-height / bestFit column-width recalculation.
# Those all stem from the file simply being re-saved by a newer Excel
# build on a different PC/user profile, not from a deliberate edit, so
# there's nothing to replay for them here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Status cell for the "Load Data" row.
$ws.Range("E2").Value = "Done"
$ws.Range("E2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E2").VerticalAlignment = -4108     # xlCenter

# Cosmetic: the saved file's last selection moved to I12.
[void]$ws.Range("I12").Select()
